# Generate Report for Handback
#
# Updates the timestamps recorded for the file
# "79e7ad8f-45e1-445a-94e1-ff455b360b30.md" (row 3 on every sheet) to
# reflect a fresh handoff/handback cycle:
#   - Overview!G3            "Latest HO Xliff Generate Date"
#   - zh-cn!H3 / zh-cn!L3     "Correspond Handoff Datetime" / "Correspond Handback DateTime"
#   - de-de!H3 / de-de!L3     "Correspond Handoff Datetime" / "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2017-01-03 05:06:17"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2017-01-03 05:06:05"
$zhcn.Range("L3").Value = "2017-01-03 05:06:38"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2017-01-03 05:06:17"
$dede.Range("L3").Value = "2017-01-03 05:06:49"
